$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '49.696.43'
$ws.Cells.Item(2, 5).Value = '  -0.76%  '

$ws.Cells.Item(3, 4).Value = '2.650.04'
$ws.Cells.Item(3, 5).Value = '  +0.07%  '

$ws.Cells.Item(4, 5).Value = '  +0.02%  '

$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '112.89'
$ws.Cells.Item(5, 5).Value = '  -1.02%  '

$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '327.39'
$ws.Cells.Item(6, 5).Value = '  +0.21%  '

$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = '0.524'
$ws.Cells.Item(7, 5).Value = '  -1.15%  '

$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = '0.999'
$ws.Cells.Item(8, 5).Value = '  -0.03%  '

$ws.Cells.Item(9, 5).Value = '  -1.12%  '

$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = '39.76'
$ws.Cells.Item(10, 5).Value = '  -3.07%  '

$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = '20.02'
$ws.Cells.Item(11, 5).Value = '  -0.65%  '

$ws.Cells.Item(12, 5).Value = '  -0.76%  '

$ws.Cells.Item(13, 5).Value = '  +2.29%  '

$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = '7.58'
$ws.Cells.Item(14, 5).Value = '  +2.64%  '

$ws.Cells.Item(15, 4).Value = '3.063.51'
$ws.Cells.Item(15, 5).Value = '  -0.01%  '

$ws.Cells.Item(16, 4).Value = '2.637.57'
$ws.Cells.Item(16, 5).Value = '  -0.93%  '

$ws.Cells.Item(17, 5).Value = '  -1.58%  '

$ws.Cells.Item(18, 4).Value = '49.676.04'

$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = '13.31'
$ws.Cells.Item(19, 5).Value = '  +0.61%  '

$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = '6.70'
$ws.Cells.Item(20, 5).Value = '  -1.32%  '

$ws.Cells.Item(21, 5).Value = '  -0.04%  '

$ws.Cells.Item(22, 5).Value = '  -0.67%  '

$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = '269.26'
$ws.Cells.Item(23, 5).Value = '  -2.54%  '

$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = '69.14'
$ws.Cells.Item(24, 5).Value = '  -4.24%  '

$ws.Cells.Item(25, 5).Value = '  -0.68%  '

$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = '26.18'
$ws.Cells.Item(26, 5).Value = '  -2.41%  '

$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = '0.999'
$ws.Cells.Item(27, 5).Value = '  -0.04%  '

$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = '10.21'
$ws.Cells.Item(28, 5).Value = '  +1.46%  '

$ws.Cells.Item(29, 5).Value = '  -0.81%  '

$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = '0.138'
$ws.Cells.Item(30, 5).Value = '  -2.29%  '

$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = '34.96'
$ws.Cells.Item(31, 5).Value = '  -3.90%  '

$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = '49.59'
$ws.Cells.Item(32, 5).Value = '  -1.43%  '

$ws.Cells.Item(33, 5).Value = '  +0.37%  '

$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = '0.0822'
$ws.Cells.Item(34, 5).Value = '  +0.84%  '

$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = '19.18'
$ws.Cells.Item(35, 5).Value = '  -2.12%  '

$ws.Cells.Item(36, 5).Value = '  -0.15%  '

$ws.Cells.Item(37, 5).Value = '  -1.83%  '

$ws.Cells.Item(38, 5).Value = '  -1.60%  '

$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = '3.13'
$ws.Cells.Item(39, 5).Value = '  +1.07%  '

$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = '129.17'
$ws.Cells.Item(40, 5).Value = '  +4.41%  '

$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = '23.63'
$ws.Cells.Item(41, 5).Value = '  +6.92%  '

$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = '0.0345'
$ws.Cells.Item(42, 5).Value = '  +9.13%  '

$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = '2.27'
$ws.Cells.Item(43, 5).Value = '  +2.22%  '

$ws.Cells.Item(44, 5).Value = '  -0.58%  '

$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = '3.34'
$ws.Cells.Item(45, 5).Value = '  -0.22%  '

$ws.Cells.Item(46, 4).Value = '2.065.88'
$ws.Cells.Item(46, 5).Value = '  -0.88%  '

$ws.Cells.Item(47, 5).Value = '  +6.24%  '

$ws.Cells.Item(48, 5).Value = '  -3.41%  '

$ws.Cells.Item(49, 5).Value = '  -2.21%  '

$ws.Cells.Item(50, 5).Value = '  -2.31%  '

$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = '58.96'
$ws.Cells.Item(51, 5).Value = '  -1.64%  '
